# Daily Orders update - data entry for 20-Jan and 21-Jan (columns AB/AC)
# Commit message: "data updated till 21 Jan 9AM"
#
# Column AB = 20-Jan, Column AC = 21-Jan.
# F/G/B2/D2/F2/G2/AB2/AC2 are formulas and recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3  -> 21-Jan (AC)
$ws.Range("AC3").Value = 3120

# Row 10 -> 20-Jan (AB)
$ws.Range("AB10").Value = 5200

# Row 14 -> 21-Jan (AC)
$ws.Range("AC14").Value = 5200

# Row 19 -> 20-Jan (AB) ; pick up the same highlight style used on the
# neighbouring 19-Jan cell (AA19) by copy/pasting it across.
$ws.Range("AA19").Copy($ws.Range("AB19"))
$ws.Range("AB19").Value = 2080

# Row 22 -> 21-Jan (AC)
$ws.Range("AC22").Value = 1040

# Row 26 -> 20-Jan (AB)
$ws.Range("AB26").Value = 2080

# Row 28 -> 21-Jan (AC)
$ws.Range("AC28").Value = 2080

# Row 31 -> 20-Jan (AB)
$ws.Range("AB31").Value = 2080

# Row 35 -> 20-Jan (AB)
$ws.Range("AB35").Value = 3120

# Row 37 -> 21-Jan (AC)
$ws.Range("AC37").Value = 3120

# Row 39 -> 20-Jan (AB)
$ws.Range("AB39").Value = 2080

# Row 44 -> 21-Jan (AC)
$ws.Range("AC44").Value = 3120

# Row 48 -> 21-Jan (AC)
$ws.Range("AC48").Value = 3120

# Row 49 -> 20-Jan (AB)
$ws.Range("AB49").Value = 3120

# Row 52 -> 21-Jan (AC)
$ws.Range("AC52").Value = 1040

# Row 53 -> 20-Jan (AB) and 21-Jan (AC)
$ws.Range("AB53").Value = 1040
$ws.Range("AC53").Value = 1040

# Row 56 -> 20-Jan (AB)
$ws.Range("AB56").Value = 2080

# Row 57 -> 21-Jan (AC)
$ws.Range("AC57").Value = 5200

# Row 59 -> 21-Jan (AC)
$ws.Range("AC59").Value = 2080

# Row 62 -> 20-Jan (AB)
$ws.Range("AB62").Value = 2080

# Row 66 -> 20-Jan (AB)
$ws.Range("AB66").Value = 3120

# Row 67 -> 21-Jan (AC)
$ws.Range("AC67").Value = 2080

# Row 71 -> 20-Jan (AB)
$ws.Range("AB71").Value = 2080

# Row 72 -> 20-Jan (AB)
$ws.Range("AB72").Value = 2080

# Row 79 -> 20-Jan (AB)
$ws.Range("AB79").Value = 3120

# Row 80 -> 20-Jan (AB)
$ws.Range("AB80").Value = 1040

# Row 81 -> 20-Jan (AB)
$ws.Range("AB81").Value = 2080

# Row 95 -> 20-Jan (AB)
$ws.Range("AB95").Value = 2080

# Leave the cursor on the last-entered cell, matching where the author's
# selection ended up after typing in the new figures.
$ws.Range("AB19").Select()
